$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Adam10-Epha3 results (rows 2-13), following Dr Hou's advice.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric statistics.
$data = @(
    @("ECs","Adam10","Epha3","FAPs",3,1,48.78557366666666,146.356721,0.3597500064892136,0.3597500064892135,3,1,46.25093466666667,138.752804,0.9569015955251317,0.9569015955251318,2256.378380332854,20307.40542299568,0.3442453551997049,0.3442453551997049),
    @("ECs","Adam10","Epha3","M2",3,1,48.78557366666666,146.356721,0.3597500064892136,0.3597500064892135,2,0.6666666666666666,0.01632333333333333,0.04897,0.0003377190931065126,0.0003377190931065127,0.7963431808188888,7.167088627369999,0.0001214944459365992,0.0001214944459365992),
    @("ECs","Adam10","Epha3","sCs",3,1,48.78557366666666,146.356721,0.3597500064892136,0.3597500064892135,3,1,2.066797333333333,6.200391999999999,0.04276068538176181,0.04276068538176181,100.8298935594035,907.4690420346318,0.01538315684357203,0.01538315684357203),
    @("FAPs","Adam10","Epha3","FAPs",3,1,22.11880033333333,66.35640099999999,0.1631063850518409,0.1631063850518409,3,1,46.25093466666667,138.752804,0.9569015955251317,0.9569015955251318,1023.015189122045,9207.136702098403,0.1560767600964431,0.1560767600964431),
    @("FAPs","Adam10","Epha3","M2",3,1,22.11880033333333,66.35640099999999,0.1631063850518409,0.1631063850518409,2,0.6666666666666666,0.01632333333333333,0.04897,0.0003377190931065126,0.0003377190931065127,0.3610525507744444,3.24947295697,0.00005508414043958937,0.00005508414043958938),
    @("FAPs","Adam10","Epha3","sCs",3,1,22.11880033333333,66.35640099999999,0.1631063850518409,0.1631063850518409,3,1,2.066797333333333,6.200391999999999,0.04276068538176181,0.04276068538176181,45.71507754546577,411.4356979091919,0.006974540814958268,0.006974540814958268),
    @("M2","Adam10","Epha3","FAPs",3,1,38.55383266666667,115.661498,0.2843000606446504,0.2843000606446503,3,1,46.25093466666667,138.752804,0.9569015955251317,0.9569015955251318,1783.150795815599,16048.35716234039,0.2720471816387576,0.2720471816387576),
    @("M2","Adam10","Epha3","M2",3,1,38.55383266666667,115.661498,0.2843000606446504,0.2843000606446503,2,0.6666666666666666,0.01632333333333333,0.04897,0.0003377190931065126,0.0003377190931065127,0.6293270618955555,5.663943557060001,0.00009601355865103785,0.00009601355865103785),
    @("M2","Adam10","Epha3","sCs",3,1,38.55383266666667,115.661498,0.2843000606446504,0.2843000606446503,3,1,2.066797333333333,6.200391999999999,0.04276068538176181,0.04276068538176181,79.68295854524622,717.1466269072159,0.0121568654472417,0.01215686544724169),
    @("sCs","Adam10","Epha3","FAPs",3,1,26.15144666666667,78.45434,0.1928435478142953,0.1928435478142952,3,1,46.25093466666667,138.752804,0.9569015955251317,0.9569015955251318,1209.528851218818,10885.75966096936,0.1845322985902262,0.1845322985902262),
    @("sCs","Adam10","Epha3","M2",3,1,26.15144666666667,78.45434,0.1928435478142953,0.1928435478142952,2,0.6666666666666666,0.01632333333333333,0.04897,0.0003377190931065126,0.0003377190931065127,0.4268787810888889,3.8419090298,0.00006512694807928621,0.00006512694807928621),
    @("sCs","Adam10","Epha3","sCs",3,1,26.15144666666667,78.45434,0.1928435478142953,0.1928435478142952,3,1,2.066797333333333,6.200391999999999,0.04276068538176181,0.04276068538176181,54.04974023347555,486.4476621012799,0.008246122275989821,0.008246122275989819)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
